$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 342; existing rows 342-375 shift down to 343-376.
$ws.Rows(342).Insert()

# Populate the newly inserted row 342 with its data.
$ws.Cells.Item(342, 1).Value = 5
$ws.Cells.Item(342, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(342, 3).Value = "Maule"
$ws.Cells.Item(342, 4).Value = 44858
$ws.Cells.Item(342, 5).Value = 7
$ws.Cells.Item(342, 6).Value = 100112003
$ws.Cells.Item(342, 7).Value = "Ajo"
$ws.Cells.Item(342, 8).Value = "Chino"
$ws.Cells.Item(342, 9).Value = "1a nueva(o)"
$ws.Cells.Item(342, 10).Value = 300
$ws.Cells.Item(342, 11).Value = 18000
$ws.Cells.Item(342, 12).Value = 18000
$ws.Cells.Item(342, 13).Value = 18000
$ws.Cells.Item(342, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(342, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(342, 16).Value = 1800
$ws.Cells.Item(342, 17).Value = 10
$ws.Cells.Item(342, 18).Value = "Hortaliza"
